# Added LOINC codes for MMC. Created stub HearingObservation.
#
# This updates the "Date" and "Count" metadata values, and collapses the
# "Concepts" sheet down to a single stub row (communicate-without-assistance),
# removing the six low/high-demand "convey/short-conversations" rows.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B8").Value = "2022-03-24T15:01:17-04:00"   # Date

# Leading "'" forces literal text "1" (matching the existing shared string
# type) instead of Excel inferring a number.
$meta.Range("B23").Value = "'1"                           # Count

# --- Concepts sheet ---------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# Drop the six low/high-demand rows (rows 3-8), keeping the header (row 1)
# and row 2 (which we overwrite below) intact - this also shrinks the used
# range/dimension down to A1:D2.
$concepts.Range("A3:D8").EntireRow.Delete()

# Replace row 2 with the single stub concept (leading "'" forces the Level
# value to be stored as literal text "1", matching the existing shared
# string used elsewhere, rather than being inferred as a number).
$concepts.Range("A2").Value = "'1"
$concepts.Range("B2").Value = "communicate-without-assistance"
$concepts.Range("C2").Value = "Participate in communication exchanges without assistance"
$concepts.Range("D2").Value = "How often does the individual participate in communication exchanges WITHOUT additional assistance from communication partner (no more than would be expected for chronological age)?"
